$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix a stray data inconsistency: G9 used the shared string "T" (capital)
# while every other row in that column ("product_usetaxes") uses lowercase "t".
$ws.Cells.Item(9, 7).Value2 = "t"

# --- New column A: mirror the existing product_code values already held in
# column B (literal values, not a formula).
for ($r = 2; $r -le 9; $r++) {
    $ws.Cells.Item($r, 1).Value2 = $ws.Cells.Item($r, 2).Value2
}

# --- Move the stray "discount" values that lived in column L into the new
# column K, and populate the rest of column K with 0 (no discount).
$ws.Range("L7").ClearContents()
$ws.Range("L8").ClearContents()

$ws.Cells.Item(1, 11).Value2 = "descuento"
$discounts = @(0, 0, 0, 0, 0, 5, 5, 0)
for ($r = 2; $r -le 9; $r++) {
    $ws.Cells.Item($r, 11).Value2 = $discounts[$r - 2]
}

# --- Highlight the whole data body (not the header row) with a yellow fill,
# matching the new "descuento" column's emphasis.
$ws.Range("A2:K9").Interior.Color = 65535

# --- Restore the cursor/selection to where the author left it.
$ws.Range("C13").Select()
